$wb = $excel.ActiveWorkbook

# "Sheet1" holds the pool of not-yet-used image ids (one per row, column A).
# "used" is the log of ids that have been consumed, with source filename + timestamp.
$names = $wb.Worksheets.Item("Sheet1")
$used  = $wb.Worksheets.Item("used")

# The two ids at the top of the pool (rows 2 and 3 - row 1 is a header-ish
# first entry that stays put) have been consumed; pop them off the pool...
$id1 = $names.Range("A2").Text
$id2 = $names.Range("A3").Text

$names.Rows.Item(2).Delete()
$names.Rows.Item(2).Delete()

# ...and append them to the "used" log with their source filenames and the
# time they were used.
$lastRow = $used.Cells.Item($used.Rows.Count, 1).End(-4162).Row
$nextRow = $lastRow + 1

$used.Range("A" + $nextRow).Value = $id1
$used.Range("B" + $nextRow).Value = "ChatGPT Image 2026年1月21日 20_56_59.png"
$used.Range("C" + $nextRow).Value = "2026-01-21 20:58:13"

$nextRow = $nextRow + 1

$used.Range("A" + $nextRow).Value = $id2
$used.Range("B" + $nextRow).Value = "ChatGPT Image 2026年1月21日 20_57_06.png"
$used.Range("C" + $nextRow).Value = "2026-01-21 20:58:13"
